$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 569; this shifts existing rows 569-604 down to 570-605
$ws.Rows.Item(569).Insert()

# Populate the new row 569 with the new record's data
$ws.Cells.Item(569, 1).Value = 10
$ws.Cells.Item(569, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(569, 3).Value = "La Araucanía"
$ws.Cells.Item(569, 4).Value = 45265
$ws.Cells.Item(569, 5).Value = 9
$ws.Cells.Item(569, 6).Value = 100112009
$ws.Cells.Item(569, 7).Value = "Acelga"
$ws.Cells.Item(569, 8).Value = "Sin especificar"
$ws.Cells.Item(569, 9).Value = "Primera"
$ws.Cells.Item(569, 10).Value = 80
$ws.Cells.Item(569, 11).Value = 10000
$ws.Cells.Item(569, 12).Value = 10000
$ws.Cells.Item(569, 13).Value = 10000
$ws.Cells.Item(569, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(569, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(569, 16).Value = 833
$ws.Cells.Item(569, 17).Value = 12
$ws.Cells.Item(569, 18).Value = "Hortaliza"
